# Fraud Detection Report: add newly-classified procedures, a totals row,
# and tidy up a couple of floating point roundings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the three new data rows, shifting existing rows down ----------
# New "Spirometry (procedure)" row lands at row 6 (between Colonoscopy and
# Standard pregnancy test).
$ws.Rows(6).Insert()
# New "Throat culture (procedure)" row lands at row 8 (between Standard
# pregnancy test and Prostatectomy).
$ws.Rows(8).Insert()
# New "Upper arm X-ray" row lands at row 9 (between Throat culture and
# Prostatectomy).
$ws.Rows(9).Insert()

# --- Minor floating point tidy-ups on pre-existing rows --------------------
$ws.Range("D3").Value = 68.59999999999999
$ws.Range("D7").Value = 43.4

# --- Row 6: Spirometry (procedure) -- Fraud (red, style of row 5) ----------
$ws.Range("A6").Value = "Spirometry (procedure)"
$ws.Range("B6").Value = 15000
$ws.Range("C6").Value = 7786.47
$ws.Range("D6").Value = 7213.53
$ws.Range("E6").Value = "Fraud"

# --- Row 8: Throat culture (procedure) -- Risk (orange, style of row 2) ---
$ws.Range("A2:E2").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Throat culture (procedure)"
$ws.Range("B8").Value = 2300
$ws.Range("C8").Value = 2020.43
$ws.Range("D8").Value = 279.57
$ws.Range("E8").Value = "Risk"

# --- Row 9: Upper arm X-ray -- Risk (orange, style of row 2) --------------
$ws.Range("A2:E2").Copy()
$ws.Range("A9:E9").PasteSpecial(-4122)
$ws.Range("A9").Value = "Upper arm X-ray"
$ws.Range("B9").Value = 1500
$ws.Range("C9").Value = 431.4
$ws.Range("D9").Value = 1068.6
$ws.Range("E9").Value = "Risk"

# --- Row 12: Total Invoice Amount summary row (new white fill) ------------
$ws.Range("A12:E12").Interior.ColorIndex = 2
$ws.Range("A12").Value = "Total Invoice Amount"
$ws.Range("B12").Value = $ws.Application.WorksheetFunction.Sum($ws.Range("B2:B11"))
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
